$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H70").Value = 1933
$ws.Range("J70").Value = 1899.5
$ws.Range("L70").Value = 5698.5
$ws.Range("N70").Value = -6238.5

$ws.Range("H73").Value = 1933
$ws.Range("J73").Value = 1899.5
$ws.Range("L73").Value = 5698.5
$ws.Range("N73").Value = -7570.5

$ws.Range("H80").Value = 683.2222
$ws.Range("J80").Value = 544.5
$ws.Range("L80").Value = 1633.5
$ws.Range("N80").Value = -3629.5

$ws.Range("H83").Value = 683.2222
$ws.Range("J83").Value = 544.5
$ws.Range("L83").Value = 4900.5
$ws.Range("N83").Value = -14884.5

$ws.Range("H111").Value = 2433.4546
$ws.Range("I111").Value = 2476.8
$ws.Range("K111").Value = 7430.400000000001
$ws.Range("M111").Value = -4363.400000000001

$ws.Range("H112").Value = 1537.7
$ws.Range("J112").Value = 1537.7
$ws.Range("L112").Value = 4613.1
$ws.Range("N112").Value = -6829.1

$ws.Range("H137").Value = 2565.6365
$ws.Range("I137").Value = 1403.4348
$ws.Range("K137").Value = 4210.3044
$ws.Range("M137").Value = -1660.3044

$ws.Range("H138").Value = 6818.8667
$ws.Range("I138").Value = 2165.5
$ws.Range("K138").Value = 6496.5
$ws.Range("M138").Value = -1356.5


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4074.8386
$ws.Range("I32").Value = 3710.6667
$ws.Range("K32").Value = 3710.6667
$ws.Range("M32").Value = -3423.6667

$ws.Range("H61").Value = 1812.3077
$ws.Range("I61").Value = 1630
$ws.Range("K61").Value = 1630
$ws.Range("M61").Value = -1418

$ws.Range("H74").Value = 15381193
$ws.Range("I74").Value = 24991314
$ws.Range("K74").Value = 24991314
$ws.Range("M74").Value = -24990440

$ws.Range("H77").Value = 15381193
$ws.Range("I77").Value = 24991314
$ws.Range("K77").Value = 124956570
$ws.Range("M77").Value = -124952202

$ws.Range("H122").Value = 1408.2727
$ws.Range("I122").Value = 1474.1
$ws.Range("J122").Value = 750
$ws.Range("K122").Value = 4422.299999999999
$ws.Range("L122").Value = 2250
$ws.Range("M122").Value = -1972.299999999999
$ws.Range("N122").Value = -7150

$ws.Range("H132").Value = 2036.0238
$ws.Range("J132").Value = 4035.0715
$ws.Range("L132").Value = 12105.2145
$ws.Range("N132").Value = -17165.2145

$ws.Range("H136").Value = 1812.3077
$ws.Range("I136").Value = 1630
$ws.Range("K136").Value = 4890
$ws.Range("M136").Value = -2340


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3605.9167
$ws.Range("I86").Value = 3503.4285
$ws.Range("J86").Value = 3749.4
$ws.Range("K86").Value = 3503.4285
$ws.Range("L86").Value = 3749.4
$ws.Range("M86").Value = -2380.4285
$ws.Range("N86").Value = -5995.4

$ws.Range("H89").Value = 3605.9167
$ws.Range("I89").Value = 3503.4285
$ws.Range("J89").Value = 3749.4
$ws.Range("K89").Value = 17517.1425
$ws.Range("L89").Value = 18747
$ws.Range("M89").Value = -11901.1425
$ws.Range("N89").Value = -29979


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4635.6665
$ws.Range("I99").Value = 4561.75
$ws.Range("J99").Value = 4783.5
$ws.Range("K99").Value = 4561.75
$ws.Range("L99").Value = 4783.5
$ws.Range("M99").Value = -3063.75
$ws.Range("N99").Value = -7779.5

$ws.Range("H114").Value = 50000.5
$ws.Range("J114").Value = 50000.5
$ws.Range("L114").Value = 50000.5
$ws.Range("N114").Value = -58678.5

$ws.Range("H122").Value = 2106.75
$ws.Range("I122").Value = 1660.25
$ws.Range("J122").Value = 2999.75
$ws.Range("K122").Value = 4980.75
$ws.Range("L122").Value = 8999.25
$ws.Range("M122").Value = -2530.75
$ws.Range("N122").Value = -13899.25

$ws.Range("H126").Value = 4635.6665
$ws.Range("I126").Value = 4561.75
$ws.Range("J126").Value = 4783.5
$ws.Range("K126").Value = 13685.25
$ws.Range("L126").Value = 14350.5
$ws.Range("M126").Value = -11215.25
$ws.Range("N126").Value = -19290.5


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 20002
$ws.Range("J9").Value = 20002
$ws.Range("L9").Value = 60006
$ws.Range("N9").Value = -60454

$ws.Range("H33").Value = 1569
$ws.Range("I33").Value = 197.2
$ws.Range("K33").Value = 1183.2
$ws.Range("M33").Value = -900.1999999999998

$ws.Range("H115").Value = 3000
$ws.Range("I115").Value = 3000
$ws.Range("K115").Value = 9000
$ws.Range("M115").Value = -7825


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()

$ws.Range("H122").Value = 1529.4445
$ws.Range("I122").Value = 1378.3334
$ws.Range("K122").Value = 4135.0002
$ws.Range("M122").Value = -1685.0002

$ws.Range("H126").Value = 2773.5
$ws.Range("I126").Value = 2550
$ws.Range("J126").Value = 2997
$ws.Range("K126").Value = 7650
$ws.Range("L126").Value = 8991
$ws.Range("M126").Value = -5180
$ws.Range("N126").Value = -13931

$ws.Range("H132").Value = 2844.6155
$ws.Range("I132").Value = 1613.8334
$ws.Range("K132").Value = 4841.5002
$ws.Range("M132").Value = -2311.5002


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 295.18182
$ws.Range("I55").Value = 293.75
$ws.Range("K55").Value = 293.75
$ws.Range("M55").Value = -120.75

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H132").Value = 3614
$ws.Range("J132").Value = 4140.6665
$ws.Range("L132").Value = 12421.9995
$ws.Range("N132").Value = -17481.9995


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11615.77
$ws.Range("I81").Value = 6374.5
$ws.Range("K81").Value = 12749
$ws.Range("M81").Value = -11688

$ws.Range("H84").Value = 11615.77
$ws.Range("I84").Value = 6374.5
$ws.Range("K84").Value = 63745
$ws.Range("M84").Value = -58441

$ws.Range("H136").Value = 2270.8572
$ws.Range("I136").Value = 1786.125
$ws.Range("J136").Value = 2917.1667
$ws.Range("K136").Value = 5358.375
$ws.Range("L136").Value = 8751.500100000001
$ws.Range("M136").Value = -2808.375
$ws.Range("N136").Value = -13851.5001

